# Applies the 'Add more sound effects, actors and fix some bugs' edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: row 10 (Id=9) category changed from 1 to 4
$ws.Range("B10").Value = 4

# New rows of questions appended at the bottom of the sheet (rows 29-71)
$newRows = @(
    ,@(29, 28, 6, 'حاجة مؤثرة جداً', 'كدة رضا - احمد حلمي', 'null')
    ,@(30, 29, 6, 'يابني ابعد عني يابني', 'الجوكر - محمد صبحي', 'null')
    ,@(31, 30, 6, 'دي حاجة لو عرفتوها تبقوا عُمد', 'محامي خلع - حسن حسني', 'null')
    ,@(32, 31, 2, 'ي س ر ك  ك ر ح ت م', 'كرسي متحرك', 'null')
    ,@(33, 32, 2, 'ة ب ع ل', 'لعبة', 'null')
    ,@(34, 33, 2, 'ة ر ك ف', 'فكرة ', 'null')
    ,@(35, 34, 2, 'س ي س م ر', 'رمسيس', 'null')
    ,@(36, 35, 2, 'ن و ع ر ف', 'فرعون', 'null')
    ,@(37, 36, 2, 'ل د ا ع  م ا م ا', 'عادل امام', 'null')
    ,@(38, 37, 2, 'ة ع ر ز م  ن ج ا و د', 'مزرعة دواجن', 'null')
    ,@(39, 38, 2, 'ق و ر ش  س م ش ل ا', 'شروق الشمس', 'null')
    ,@(40, 39, 2, ' ل ا ي خ', 'خيال', 'null')
    ,@(41, 40, 2, 'ل و ف  س م د م', 'فول مدمس', 'null')
    ,@(42, 41, 6, 'انا بابا يلا', 'تيتو - خالد صالح', 'null')
    ,@(43, 42, 6, 'هنرقص دانص يا روح امك', 'الفرن - عادل ادهم', 'null')
    ,@(44, 43, 6, 'جبت الترنك', 'الناظر - حسين ابو حجاج', 'null')
    ,@(45, 44, 6, 'جي تعزي ولا جي تهزر', 'الكيف - يحيي الفخراني', 'null')
    ,@(46, 45, 7, 'ايه الصوت ده؟', 'Popcorn', 'Resources/Soundeffects/popcorn.mp3')
    ,@(47, 46, 7, 'ايه الصوت ده؟', 'Pouring Tea', 'Resources/Soundeffects/tea.mp3')
    ,@(48, 47, 7, 'ايه الصوت ده؟', 'Thunder', 'Resources/Soundeffects/thunder.mp3')
    ,@(49, 48, 1, '5 مسلسلات من رمضان 2024', 'بابا جه - الحشاشين - امبراطورية م - عتبات البهجة - مسار اجباري', 'null')
    ,@(50, 49, 1, 'Smallest Continent in Area', 'Australia', 'null')
    ,@(51, 50, 1, 'Largest Continet in Area', 'Asia', 'null')
    ,@(52, 51, 1, 'Country that contains 7000 island', 'Phillipins ', 'null')
    ,@(53, 52, 1, 'The only Arabic country with no desert', 'Lebanon', 'null')
    ,@(54, 53, 1, 'Who had drawn the Mona Lisa', 'Leonardo Davinci ', 'null')
    ,@(55, 54, 1, 'حاجة بتتاكل عمرها ما بتبوظ', 'العسل', 'null')
    ,@(56, 55, 1, 'اكبر كوكب في المجرة الشمسية؟', 'المشترى - Jupiter ', 'null')
    ,@(57, 56, 1, 'المغرب في انهي قارة؟', 'افريقيا', 'null')
    ,@(58, 57, 8, 'مين ده؟', 'Tom Hanks', 'Resources/Actors/TomHanks.jpg')
    ,@(59, 58, 1, 'اكبر عضو في جسم الإنسان؟', 'الجلد', 'null')
    ,@(60, 59, 1, 'علم البرازيل معمول من كام لون؟', '4 (Yellow, White, Green, Blue)', 'null')
    ,@(61, 60, 6, 'امال لأ امال طبعاً', 'حزمني يا - شريف منير', 'null')
    ,@(62, 61, 6, 'دي اراء ارااااء', 'رمضان مبروك ابو العلمين حمودة - محمد هنيدي', 'null')
    ,@(63, 62, 6, 'انا مش تبع حد انا رئيس جمهورية نفسي', 'ظرف طارق - محمد شرف', 'null')
    ,@(64, 63, 4, 'اغنية لـ حمادة هلال', 'null', 'null')
    ,@(65, 64, 4, 'اغنية لـ محمد هنيدي', 'null', 'null')
    ,@(66, 65, 4, 'اغنية لـ أصالة', 'null', 'null')
    ,@(67, 66, 4, 'اغنية لـ سميرة سعيد', 'null', 'null')
    ,@(68, 67, 4, 'اغنية للشاب خالد', 'null', 'null')
    ,@(69, 68, 4, 'اغنية لـ شرين ', 'null', 'null')
    ,@(70, 69, 4, 'اغنية لـ MTM', 'null', 'null')
    ,@(71, 70, 4, 'اغنية لـ وسط البلد', 'null', 'null')
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# Autofit columns C and D to match new, wider content
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Move selection / viewport to the end of the data, like the author's last edit
$ws.Range("E71").Select()
$excel.ActiveWindow.ScrollRow = 59
